$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = -2.796573638916016
$ws.Cells.Item(2, 2).Value = -10.2697286605835
$ws.Cells.Item(2, 3).Value = 2.567607164382935
$ws.Cells.Item(3, 1).Value = 0.8655490875244141
$ws.Cells.Item(3, 2).Value = -18.99405670166016
$ws.Cells.Item(3, 3).Value = 5.860441207885742
$ws.Cells.Item(4, 1).Value = 12.034010887146
$ws.Cells.Item(4, 2).Value = -13.59659099578857
$ws.Cells.Item(4, 3).Value = 12.0263729095459
$ws.Cells.Item(5, 1).Value = 5.544707775115967
$ws.Cells.Item(5, 2).Value = -11.79852485656738
$ws.Cells.Item(5, 3).Value = -1.972949981689453
$ws.Cells.Item(6, 1).Value = 6.281956195831299
$ws.Cells.Item(6, 2).Value = -19.22240829467773
$ws.Cells.Item(6, 3).Value = -3.653380393981934
$ws.Cells.Item(7, 1).Value = -1.975251197814941
$ws.Cells.Item(7, 2).Value = -43.3292350769043
$ws.Cells.Item(7, 3).Value = 0.271059513092041
$ws.Cells.Item(8, 1).Value = -20.38084411621094
$ws.Cells.Item(8, 2).Value = -15.68858814239502
$ws.Cells.Item(8, 3).Value = 1.961378574371338
$ws.Cells.Item(9, 1).Value = 10.76584243774414
$ws.Cells.Item(9, 2).Value = -33.70746231079102
$ws.Cells.Item(9, 3).Value = 31.38712692260743
$ws.Cells.Item(10, 1).Value = -27.7362232208252
$ws.Cells.Item(10, 2).Value = -0.916855812072754
$ws.Cells.Item(10, 3).Value = -10.06494903564453
$ws.Cells.Item(11, 1).Value = 8.33868408203125
$ws.Cells.Item(11, 2).Value = -7.790350914001465
$ws.Cells.Item(11, 3).Value = 11.65683746337891
$ws.Cells.Item(12, 1).Value = -1.120648384094239
$ws.Cells.Item(12, 2).Value = -23.34181785583496
$ws.Cells.Item(12, 3).Value = -21.7900619506836
$ws.Cells.Item(13, 1).Value = 33.84098815917969
$ws.Cells.Item(13, 2).Value = -33.77373504638672
$ws.Cells.Item(13, 3).Value = -8.738304138183594
$ws.Cells.Item(14, 1).Value = -9.054259300231934
$ws.Cells.Item(14, 2).Value = -1.461036801338196
$ws.Cells.Item(14, 3).Value = 6.970683097839356
$ws.Cells.Item(15, 1).Value = 10.75043201446533
$ws.Cells.Item(15, 2).Value = -19.06211471557617
$ws.Cells.Item(15, 3).Value = 15.70715045928955
$ws.Cells.Item(16, 1).Value = 1.280778884887695
$ws.Cells.Item(16, 2).Value = -11.42607116699219
$ws.Cells.Item(16, 3).Value = -8.684724807739258
$ws.Cells.Item(17, 1).Value = 36.76531219482422
$ws.Cells.Item(17, 2).Value = -8.253963470458984
$ws.Cells.Item(17, 3).Value = -22.78386306762696
$ws.Cells.Item(18, 1).Value = -36.61545944213867
$ws.Cells.Item(18, 2).Value = -18.09431838989257
$ws.Cells.Item(18, 3).Value = 3.823569297790527
$ws.Cells.Item(19, 1).Value = 18.00795745849609
$ws.Cells.Item(19, 2).Value = -45.09830856323242
$ws.Cells.Item(19, 3).Value = -8.873518943786621
$ws.Cells.Item(20, 1).Value = -20.05809783935547
$ws.Cells.Item(20, 2).Value = 0.3998336791992187
$ws.Cells.Item(20, 3).Value = 3.845695495605469
$ws.Cells.Item(21, 1).Value = 7.206372261047363
$ws.Cells.Item(21, 2).Value = -20.33248519897461
$ws.Cells.Item(21, 3).Value = 22.94344902038575
$ws.Cells.Item(22, 1).Value = -10.1914873123169
$ws.Cells.Item(22, 2).Value = -12.15236282348633
$ws.Cells.Item(22, 3).Value = -10.96279335021973
$ws.Cells.Item(23, 1).Value = 28.82599258422852
$ws.Cells.Item(23, 2).Value = 6.167891502380371
$ws.Cells.Item(23, 3).Value = -0.517308235168457
$ws.Cells.Item(24, 1).Value = -22.66286087036133
$ws.Cells.Item(24, 2).Value = -15.7267017364502
$ws.Cells.Item(24, 3).Value = 0.2342269420623779
$ws.Cells.Item(25, 1).Value = -11.30067539215088
$ws.Cells.Item(25, 2).Value = -54.94432067871094
$ws.Cells.Item(25, 3).Value = 17.55831718444824
$ws.Cells.Item(26, 1).Value = -17.29559326171875
$ws.Cells.Item(26, 2).Value = 4.657787322998047
$ws.Cells.Item(26, 3).Value = -2.186375617980957
$ws.Cells.Item(27, 1).Value = 20.14034080505371
$ws.Cells.Item(27, 2).Value = -19.13811683654785
$ws.Cells.Item(27, 3).Value = 23.86569976806641
$ws.Cells.Item(28, 1).Value = 13.2857141494751
$ws.Cells.Item(28, 2).Value = -19.34296989440918
$ws.Cells.Item(28, 3).Value = -3.264841318130493
$ws.Cells.Item(29, 1).Value = -11.63338565826416
$ws.Cells.Item(29, 2).Value = 5.394529819488525
$ws.Cells.Item(29, 3).Value = -0.188831090927124
$ws.Cells.Item(30, 1).Value = -38.69764709472656
$ws.Cells.Item(30, 2).Value = -13.45611763000488
$ws.Cells.Item(30, 3).Value = 1.173340797424316
$ws.Cells.Item(31, 1).Value = 9.269144058227541
$ws.Cells.Item(31, 2).Value = -44.79425811767578
$ws.Cells.Item(31, 3).Value = 1.448012948036194
